$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aoki 2005")
$ws.Range("E7").NumberFormat = "0.000"
$ws.Range("E7").Value = 0.81
Write-Output ($ws.Range("E7").NumberFormat)
